# generacion de formato constancia
#
# Applies:
#  1. Replace the "No. «CORRELATIVO»-2026" MERGEFIELD construct in the
#     header with the literal placeholder text "No. ${fecha}" (dropping
#     the yellow highlight formatting that lived on the field).
#  2. Remove the stray _GoBack bookmark left over from the last cursor
#     position.
#  3. Re-key the signature block ("Licda. María José Samayoa Aldana" /
#     "Directora de Desarrollo Social" / "Municipalidad de Guatemala")
#     into the multi-run + proofErr (spell-check) shape Word produces
#     when the text is retyped through the UI.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "No. «CORRELATIVO»-2026"  ->  "No. ${fecha}"
# ---------------------------------------------------------------------
# Delete the MERGEFIELD CORRELATIVO field outright (this removes the
# begin/instrText/separate/result/end run cluster cleanly), leaving
# behind "No. " followed by the literal "-2026" runs.
$correlativo = $d.Fields(1)
$correlativo.Delete()

$noPara = $d.Paragraphs(2)
$tail = $d.Range($noPara.Range.Start + 4, $noPara.Range.End - 1)
$fechaFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:b/><w:bCs/><w:lang w:val="es-GT"/></w:rPr><w:t>${fecha}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tail.InsertXML($fechaFrag)

# ---------------------------------------------------------------------
# 2) Drop the leftover _GoBack bookmark
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 3) Signature block re-key (proofErr-wrapped runs)
# ---------------------------------------------------------------------
$rprPtBr = '<w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:b/><w:bCs/><w:szCs w:val="24"/><w:lang w:val="pt-BR"/>'

function New-WordXmlFragment($innerBody) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerBody + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# 3a) "Licda. María José Samayoa Aldana" (paragraph already has "Lic" in
#     its own run; only the remainder "da. María José Samayoa Aldana" is
#     rewritten).
$namePara = $d.Paragraphs(20)
$nameTail = $d.Range($namePara.Range.Start + 3, $namePara.Range.End - 1)
$nameInner = (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t>da</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t xml:space="preserve">. </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t>María</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t xml:space="preserve"> José </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t>Samayoa</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t>Aldana</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)
$nameTail.InsertXML((New-WordXmlFragment $nameInner))

# 3b) "Directora de Desarrollo Social"
$titlePara = $d.Paragraphs(21)
$titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
$titleInner = (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t>Directora</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t xml:space="preserve"> de </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t>Desarrollo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t xml:space="preserve"> Social</w:t></w:r>'
)
$titleRange.InsertXML((New-WordXmlFragment $titleInner))

# 3c) "Municipalidad de Guatemala"
$muniPara = $d.Paragraphs(22)
$muniRange = $d.Range($muniPara.Range.Start, $muniPara.Range.End - 1)
$muniInner = (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t>Municipalidad</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr>' + $rprPtBr + '</w:rPr><w:t xml:space="preserve"> de Guatemala</w:t></w:r>'
)
$muniRange.InsertXML((New-WordXmlFragment $muniInner))

Write-Output "done"
